{"js": "// Add \"Major \" to the front of the \"Responsibilities:\" line for each of the\n// three team members, plus member-specific trailing additions:\n//   - Ray:    append \", Readme\" after \"...walkthrough\"\n//   - Jason:  append a trailing space after \"...Adding items to game\"\n//   - Hunter: no trailing addition (only the \"Major \" prefix)\n//\n// Matching is done on the paragraph's full current text so the edit is only\n// ever applied to the exact paragraphs touched by the source change (and is\n// a no-op / safely skipped if the document doesn't contain them).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst edits = [\n  {\n    match: \"Responsibilities: General Troubleshooting/Debugging, Adding items to the game world, walkthrough\",\n    prefix: \"Major \",\n    suffix: \", Readme\",\n  },\n  {\n    match: \"Responsibilities: Bookkeeping items, Self-Assessment document, Game planning, Adding items to game\",\n    prefix: \"Major \",\n    suffix: \" \",\n  },\n  {\n    match: \"Responsibilities: General Troubleshooting/Debugging, Walkthrough, Game World Creation\",\n    prefix: \"Major \",\n    suffix: null,\n  },\n];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  const edit = edits.find((e) => text === e.match);\n  if (!edit) continue;\n\n  if (edit.prefix) {\n    para.insertText(edit.prefix, \"Start\");\n  }\n  if (edit.suffix) {\n    para.insertText(edit.suffix, \"End\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Add \"Major \" to the front of the \"Responsibilities:\" line for each of the\n# three team members, plus member-specific trailing additions:\n#   - Ray:    append \", Readme\" after \"...walkthrough\"\n#   - Jason:  append a trailing space after \"...Adding items to game\"\n#   - Hunter: no trailing addition (only the \"Major \" prefix)\n#\n# Matching is done on each paragraph's current text (end-of-paragraph mark\n# trimmed off) so the edit only ever touches the exact paragraphs affected\n# by the source change.\n\n$d = $word.ActiveDocument\n\n$edits = @(\n    @{\n        Match  = \"Responsibilities: General Troubleshooting/Debugging, Adding items to the game world, walkthrough\"\n        Prefix = \"Major \"\n        Suffix = \", Readme\"\n    },\n    @{\n        Match  = \"Responsibilities: Bookkeeping items, Self-Assessment document, Game planning, Adding items to game\"\n        Prefix = \"Major \"\n        Suffix = \" \"\n    },\n    @{\n        Match  = \"Responsibilities: General Troubleshooting/Debugging, Walkthrough, Game World Creation\"\n        Prefix = \"Major \"\n        Suffix = $null\n    }\n)\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n\n    foreach ($edit in $edits) {\n        if ($text -eq $edit.Match) {\n            if ($edit.Prefix) {\n                $p.Range.InsertBefore($edit.Prefix)\n            }\n            if ($edit.Suffix) {\n                $p.Range.InsertAfter($edit.Suffix)\n            }\n            break\n        }\n    }\n}\n"}
